$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 253.5625  # H33
$ws.Cells.Item(33, 10).Value = 306.4  # J33
$ws.Cells.Item(33, 12).Value = 306.4  # L33
$ws.Cells.Item(33, 14).Value = -764.4  # N33
$ws.Cells.Item(64, 8).Value = 4746.5  # H64
$ws.Cells.Item(64, 10).Value = 4746.5  # J64
$ws.Cells.Item(64, 12).Value = 4746.5  # L64
$ws.Cells.Item(64, 14).Value = -5242.5  # N64
$ws.Cells.Item(67, 8).Value = 4746.5  # H67
$ws.Cells.Item(67, 10).Value = 4746.5  # J67
$ws.Cells.Item(67, 12).Value = 4746.5  # L67
$ws.Cells.Item(67, 14).Value = -6462.5  # N67
$ws.Cells.Item(70, 8).Value = 8053.5  # H70
$ws.Cells.Item(70, 9).Value = 7143.3335  # I70
$ws.Cells.Item(70, 10).Value = 8736.125  # J70
$ws.Cells.Item(70, 11).Value = 21430.0005  # K70
$ws.Cells.Item(70, 12).Value = 26208.375  # L70
$ws.Cells.Item(70, 13).Value = -21160.0005  # M70
$ws.Cells.Item(70, 14).Value = -26748.375  # N70
$ws.Cells.Item(73, 8).Value = 8053.5  # H73
$ws.Cells.Item(73, 9).Value = 7143.3335  # I73
$ws.Cells.Item(73, 10).Value = 8736.125  # J73
$ws.Cells.Item(73, 11).Value = 21430.0005  # K73
$ws.Cells.Item(73, 12).Value = 26208.375  # L73
$ws.Cells.Item(73, 13).Value = -20494.0005  # M73
$ws.Cells.Item(73, 14).Value = -28080.375  # N73
$ws.Cells.Item(76, 8).Value = 5499.846  # H76
$ws.Cells.Item(76, 10).Value = 0  # J76
$ws.Cells.Item(76, 12).Value = 0  # L76
$ws.Cells.Item(76, 14).ClearContents()  # N76
$ws.Cells.Item(79, 8).Value = 5499.846  # H79
$ws.Cells.Item(79, 10).Value = 0  # J79
$ws.Cells.Item(79, 12).Value = 0  # L79
$ws.Cells.Item(79, 14).ClearContents()  # N79
$ws.Cells.Item(131, 8).Value = 10432.5  # H131
$ws.Cells.Item(131, 9).Value = 1079.2858  # I131
$ws.Cells.Item(131, 10).Value = 19785.715  # J131
$ws.Cells.Item(131, 11).Value = 3237.8574  # K131
$ws.Cells.Item(131, 12).Value = 59357.145  # L131
$ws.Cells.Item(131, 13).Value = 1802.1426  # M131
$ws.Cells.Item(131, 14).Value = -69437.145  # N131
$ws.Cells.Item(133, 8).Value = 86999  # H133
$ws.Cells.Item(133, 10).Value = 86999  # J133
$ws.Cells.Item(133, 12).Value = 86999  # L133
$ws.Cells.Item(133, 14).Value = -97119  # N133
$ws.Cells.Item(138, 8).Value = 3607.5417  # H138
$ws.Cells.Item(138, 9).Value = 3689.1333  # I138
$ws.Cells.Item(138, 10).Value = 3586.07  # J138
$ws.Cells.Item(138, 11).Value = 11067.3999  # K138
$ws.Cells.Item(138, 12).Value = 10758.21  # L138
$ws.Cells.Item(138, 13).Value = -5927.3999  # M138
$ws.Cells.Item(138, 14).Value = -21038.21  # N138

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4708.1665  # H32
$ws.Cells.Item(32, 9).Value = 4661.2646  # I32
$ws.Cells.Item(32, 10).Value = 5505.5  # J32
$ws.Cells.Item(32, 11).Value = 4661.2646  # K32
$ws.Cells.Item(32, 12).Value = 5505.5  # L32
$ws.Cells.Item(32, 13).Value = -4374.2646  # M32
$ws.Cells.Item(32, 14).Value = -6079.5  # N32
$ws.Cells.Item(38, 8).Value = 69614.5  # H38
$ws.Cells.Item(38, 9).Value = 69614.5  # I38
$ws.Cells.Item(38, 11).Value = 69614.5  # K38
$ws.Cells.Item(38, 13).Value = -69147.5  # M38
$ws.Cells.Item(45, 8).Value = 2070.2856  # H45
$ws.Cells.Item(45, 9).Value = 1998.8  # I45
$ws.Cells.Item(45, 11).Value = 1998.8  # K45
$ws.Cells.Item(45, 13).Value = -1621.8  # M45
$ws.Cells.Item(55, 8).Value = 19348.666  # H55
$ws.Cells.Item(55, 9).Value = 15023  # I55
$ws.Cells.Item(55, 10).Value = 28000  # J55
$ws.Cells.Item(55, 11).Value = 15023  # K55
$ws.Cells.Item(55, 12).Value = 28000  # L55
$ws.Cells.Item(55, 13).Value = -14708  # M55
$ws.Cells.Item(55, 14).Value = -28630  # N55
$ws.Cells.Item(61, 8).Value = 5333.125  # H61
$ws.Cells.Item(61, 9).Value = 1905.7142  # I61
$ws.Cells.Item(61, 11).Value = 1905.7142  # K61
$ws.Cells.Item(61, 13).Value = -1693.7142  # M61
$ws.Cells.Item(103, 8).Value = 81332.336  # H103
$ws.Cells.Item(103, 10).Value = 81332.336  # J103
$ws.Cells.Item(103, 12).Value = 81332.336  # L103
$ws.Cells.Item(103, 14).Value = -83676.336  # N103
$ws.Cells.Item(128, 8).Value = 45999.5  # H128
$ws.Cells.Item(128, 10).Value = 45999.5  # J128
$ws.Cells.Item(128, 12).Value = 45999.5  # L128
$ws.Cells.Item(128, 14).Value = -55959.5  # N128
$ws.Cells.Item(132, 8).Value = 1803.1578  # H132
$ws.Cells.Item(132, 9).Value = 1624.7778  # I132
$ws.Cells.Item(132, 11).Value = 4874.3334  # K132
$ws.Cells.Item(132, 13).Value = -2344.3334  # M132
$ws.Cells.Item(136, 8).Value = 5333.125  # H136
$ws.Cells.Item(136, 9).Value = 1905.7142  # I136
$ws.Cells.Item(136, 11).Value = 5717.142599999999  # K136
$ws.Cells.Item(136, 13).Value = -3167.142599999999  # M136

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2877.3  # H86
$ws.Cells.Item(86, 9).Value = 2515.375  # I86
$ws.Cells.Item(86, 10).Value = 4325  # J86
$ws.Cells.Item(86, 11).Value = 2515.375  # K86
$ws.Cells.Item(86, 12).Value = 4325  # L86
$ws.Cells.Item(86, 13).Value = -1392.375  # M86
$ws.Cells.Item(86, 14).Value = -6571  # N86
$ws.Cells.Item(89, 8).Value = 2877.3  # H89
$ws.Cells.Item(89, 9).Value = 2515.375  # I89
$ws.Cells.Item(89, 10).Value = 4325  # J89
$ws.Cells.Item(89, 11).Value = 12576.875  # K89
$ws.Cells.Item(89, 12).Value = 21625  # L89
$ws.Cells.Item(89, 13).Value = -6960.875  # M89
$ws.Cells.Item(89, 14).Value = -32857  # N89
$ws.Cells.Item(94, 8).Value = 27500  # H94
$ws.Cells.Item(94, 9).Value = 27500  # I94
$ws.Cells.Item(94, 11).Value = 27500  # K94
$ws.Cells.Item(94, 13).Value = -27049  # M94
$ws.Cells.Item(105, 8).Value = 3306.4  # H105
$ws.Cells.Item(105, 9).Value = 2584.2856  # I105
$ws.Cells.Item(105, 11).Value = 2584.2856  # K105
$ws.Cells.Item(105, 13).Value = -837.2856000000002  # M105
$ws.Cells.Item(107, 8).Value = 2979.4119  # H107
$ws.Cells.Item(107, 9).Value = 2664.8125  # I107
$ws.Cells.Item(107, 11).Value = 2664.8125  # K107
$ws.Cells.Item(107, 13).Value = -744.8125  # M107

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 33869.625  # H22
$ws.Cells.Item(22, 9).Value = 3488.25  # I22
$ws.Cells.Item(22, 11).Value = 3488.25  # K22
$ws.Cells.Item(22, 13).Value = -3138.25  # M22
$ws.Cells.Item(31, 8).Value = 3391.5  # H31
$ws.Cells.Item(31, 9).Value = 3159.125  # I31
$ws.Cells.Item(31, 11).Value = 3159.125  # K31
$ws.Cells.Item(31, 13).Value = -2864.125  # M31
$ws.Cells.Item(34, 8).Value = 3391.5  # H34
$ws.Cells.Item(34, 9).Value = 3159.125  # I34
$ws.Cells.Item(34, 11).Value = 3159.125  # K34
$ws.Cells.Item(34, 13).Value = -2957.125  # M34
$ws.Cells.Item(74, 8).Value = 0  # H74
$ws.Cells.Item(74, 10).Value = 0  # J74
$ws.Cells.Item(74, 12).ClearContents()  # L74
$ws.Cells.Item(74, 14).Value = 0  # N74
$ws.Cells.Item(77, 8).Value = 0  # H77
$ws.Cells.Item(77, 10).Value = 0  # J77
$ws.Cells.Item(77, 12).ClearContents()  # L77
$ws.Cells.Item(77, 14).Value = 0  # N77
$ws.Cells.Item(86, 8).Value = 6895.5  # H86
$ws.Cells.Item(86, 9).Value = 8526  # I86
$ws.Cells.Item(86, 10).Value = 5265  # J86
$ws.Cells.Item(86, 11).Value = 8526  # K86
$ws.Cells.Item(86, 12).Value = 5265  # L86
$ws.Cells.Item(86, 13).Value = -7403  # M86
$ws.Cells.Item(86, 14).Value = -7511  # N86
$ws.Cells.Item(89, 8).Value = 6895.5  # H89
$ws.Cells.Item(89, 9).Value = 8526  # I89
$ws.Cells.Item(89, 10).Value = 5265  # J89
$ws.Cells.Item(89, 11).Value = 42630  # K89
$ws.Cells.Item(89, 12).Value = 26325  # L89
$ws.Cells.Item(89, 13).Value = -37014  # M89
$ws.Cells.Item(89, 14).Value = -37557  # N89
$ws.Cells.Item(94, 8).Value = 163462.14  # H94
$ws.Cells.Item(94, 9).Value = 551512  # I94
$ws.Cells.Item(94, 11).Value = 551512  # K94
$ws.Cells.Item(94, 13).Value = -551061  # M94
$ws.Cells.Item(105, 8).Value = 3908.2222  # H105
$ws.Cells.Item(105, 9).Value = 3064.2856  # I105
$ws.Cells.Item(105, 11).Value = 3064.2856  # K105
$ws.Cells.Item(105, 13).Value = -1317.2856  # M105
$ws.Cells.Item(132, 8).Value = 1428.5  # H132
$ws.Cells.Item(132, 9).Value = 1201.4615  # I132
$ws.Cells.Item(132, 11).Value = 3604.3845  # K132
$ws.Cells.Item(132, 13).Value = -1074.3845  # M132
$ws.Cells.Item(141, 8).Value = 355496.12  # H141
$ws.Cells.Item(141, 10).Value = 355496.12  # J141
$ws.Cells.Item(141, 12).Value = 355496.12  # L141
$ws.Cells.Item(141, 14).Value = -365856.12  # N141

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 1497.4  # H14
$ws.Cells.Item(14, 9).Value = 1497.4  # I14
$ws.Cells.Item(14, 11).Value = 4492.200000000001  # K14
$ws.Cells.Item(14, 13).Value = -4319.200000000001  # M14
$ws.Cells.Item(51, 8).Value = 2546.6667  # H51
$ws.Cells.Item(51, 9).Value = 2597.3333  # I51
$ws.Cells.Item(51, 11).Value = 7791.999899999999  # K51
$ws.Cells.Item(51, 13).Value = -7331.999899999999  # M51
$ws.Cells.Item(58, 8).Value = 0  # H58
$ws.Cells.Item(58, 10).Value = 0  # J58
$ws.Cells.Item(58, 12).ClearContents()  # L58
$ws.Cells.Item(58, 14).Value = 0  # N58
$ws.Cells.Item(122, 8).Value = 506.14285  # H122
$ws.Cells.Item(122, 10).Value = 502.4  # J122
$ws.Cells.Item(122, 12).Value = 4521.599999999999  # L122
$ws.Cells.Item(122, 14).Value = -9421.599999999999  # N122
$ws.Cells.Item(131, 8).Value = 2879  # H131
$ws.Cells.Item(131, 10).Value = 3099.1667  # J131
$ws.Cells.Item(131, 12).Value = 9297.500100000001  # L131
$ws.Cells.Item(131, 14).Value = -19377.5001  # N131
$ws.Cells.Item(134, 8).Value = 16216  # H134
$ws.Cells.Item(134, 9).Value = 1000  # I134
$ws.Cells.Item(134, 11).Value = 3000  # K134
$ws.Cells.Item(134, 13).Value = 2070  # M134

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(4, 8).Value = 35683.668  # H4
$ws.Cells.Item(4, 9).Value = 42239.8  # I4
$ws.Cells.Item(4, 11).Value = 42239.8  # K4
$ws.Cells.Item(4, 13).Value = -42127.8  # M4
$ws.Cells.Item(35, 8).Value = 8000  # H35
$ws.Cells.Item(35, 9).Value = 8000  # I35
$ws.Cells.Item(35, 11).Value = 8000  # K35
$ws.Cells.Item(35, 13).Value = -7702  # M35
$ws.Cells.Item(113, 8).Value = 1437.909  # H113
$ws.Cells.Item(113, 9).Value = 1437.909  # I113
$ws.Cells.Item(113, 11).Value = 1437.909  # K113
$ws.Cells.Item(113, 13).Value = 732.0909999999999  # M113
$ws.Cells.Item(132, 8).Value = 1697.6666  # H132
$ws.Cells.Item(132, 9).Value = 1697.6666  # I132
$ws.Cells.Item(132, 11).Value = 5092.9998  # K132
$ws.Cells.Item(132, 13).Value = -2562.9998  # M132

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 507.5  # H9
$ws.Cells.Item(9, 9).Value = 410  # I9
$ws.Cells.Item(9, 10).Value = 800  # J9
$ws.Cells.Item(9, 11).Value = 410  # K9
$ws.Cells.Item(9, 12).Value = 800  # L9
$ws.Cells.Item(9, 13).Value = -186  # M9
$ws.Cells.Item(9, 14).Value = -1248  # N9
$ws.Cells.Item(35, 8).Value = 1540.5  # H35
$ws.Cells.Item(35, 9).Value = 31  # I35
$ws.Cells.Item(35, 10).Value = 3050  # J35
$ws.Cells.Item(35, 11).Value = 31  # K35
$ws.Cells.Item(35, 12).Value = 3050  # L35
$ws.Cells.Item(35, 13).Value = 305  # M35
$ws.Cells.Item(35, 14).Value = -3722  # N35
$ws.Cells.Item(40, 8).Value = 3848.8572  # H40
$ws.Cells.Item(40, 9).Value = 2869.2856  # I40
$ws.Cells.Item(40, 11).Value = 2869.2856  # K40
$ws.Cells.Item(40, 13).Value = -2733.2856  # M40
$ws.Cells.Item(100, 8).Value = 3593.2222  # H100
$ws.Cells.Item(100, 9).Value = 1359.75  # I100
$ws.Cells.Item(100, 10).Value = 5380  # J100
$ws.Cells.Item(100, 11).Value = 1359.75  # K100
$ws.Cells.Item(100, 12).Value = 5380  # L100
$ws.Cells.Item(100, 13).Value = -818.75  # M100
$ws.Cells.Item(100, 14).Value = -6462  # N100

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 15814.333  # H41
$ws.Cells.Item(41, 9).Value = 13945  # I41
$ws.Cells.Item(41, 11).Value = 13945  # K41
$ws.Cells.Item(41, 13).Value = -13555  # M41
